$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet "PartOfBodylinkTest" as the last (3rd) tab.
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PartOfBodylinkTest"

# ---------------------------------------------------------------------
# 2. Populate PartOfBodylinkTest.
#    Cells are written in the same order the strings were first
#    authored so newly-created shared strings line up the same way.
# ---------------------------------------------------------------------

# Headers / test-id columns (reuse existing shared strings).
$ws3.Range("A1").Value = "TestID"
$ws3.Range("B1").Value = "Description"
$ws3.Range("A2").Value = "TC01"

# Row 2 "Head - ..." values (new shared strings), authored in this order.
$ws3.Range("F2").Value = "Head - Ear(s)"
$ws3.Range("H2").Value = "Head - Nose"
$ws3.Range("I2").Value = "Head - Teeth"
$ws3.Range("J2").Value = "Head - Mouth"
$ws3.Range("K2").Value = "Head - Soft Tissue"
$ws3.Range("L2").Value = "Head - Facial Bones"
$ws3.Range("M2").Value = "Head - Multiple Neck Injury"
$ws3.Range("N2").Value = "Head - Vertebrae"
$ws3.Range("C2").Value = "Head - Multiple Head Injury"
$ws3.Range("D2").Value = "Head - Skull"
$ws3.Range("E2").Value = "Head - Brain"
$ws3.Range("B2").Value = "Check whether the given values are available in Part of Body table"

# Row 1 "Item N" headers (new shared strings).
$ws3.Range("C1").Value = "Item 1"
$ws3.Range("D1").Value = "Item 2"
$ws3.Range("E1").Value = "Item 3"
$ws3.Range("F1").Value = "Item 4"
$ws3.Range("G1").Value = "Item 5"
$ws3.Range("H1").Value = "Item 6"
$ws3.Range("I1").Value = "Item 7"
$ws3.Range("J1").Value = "Item 8"
$ws3.Range("K1").Value = "Item 9"
$ws3.Range("L1").Value = "Item 10"
$ws3.Range("M1").Value = "Item 11"
$ws3.Range("N1").Value = "Item 12"

# Reuse of an existing shared string ("Head - Eyes").
$ws3.Range("G2").Value = "Head - Eyes"

# ---------------------------------------------------------------------
# 3. Formatting for PartOfBodylinkTest.
# ---------------------------------------------------------------------

# Header row (A1:B1) - centered / wrapped bold-ish style matching the
# other sheets' header formatting.
$headerRange = $ws3.Range("A1:B1")
$headerRange.WrapText = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# Test id / description cells (A2:B2) - wrapped text style.
$ws3.Range("A2:B2").WrapText = $true

# Remaining data cells - plain wrap-text style.
$ws3.Range("C1:N2").WrapText = $true

# Row 2 is tall to accommodate the wrapped description text.
$ws3.Rows.Item(2).RowHeight = 67.5

# Column widths (approximate best-fit widths for the content).
$ws3.Columns.Item(2).ColumnWidth = 14.877604166666666
$ws3.Columns.Item(3).ColumnWidth = 25.592447916666668
$ws3.Columns.Item(4).ColumnWidth = 10.592447916666666
$ws3.Columns.Item(5).ColumnWidth = 10.877604166666666
$ws3.Columns.Item(6).ColumnWidth = 11.307291666666666
$ws3.Columns.Item(7).ColumnWidth = 10.307291666666666
$ws3.Columns.Item(8).ColumnWidth = 10.877604166666666
$ws3.Columns.Item(9).ColumnWidth = 11.451822916666666
$ws3.Columns.Item(10).ColumnWidth = 12.307291666666666
$ws3.Columns.Item(11).ColumnWidth = 16.022135416666668
$ws3.Columns.Item(12).ColumnWidth = 17.451822916666668
$ws3.Columns.Item(13).ColumnWidth = 24.877604166666668
$ws3.Columns.Item(14).ColumnWidth = 15.451822916666666

$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. Selections left on each sheet (as recorded by each sheet's view).
# ---------------------------------------------------------------------
[void]$sheet1.Activate()
[void]$sheet1.Range("B2").Select()

[void]$sheet2.Activate()
[void]$sheet2.Range("B3").Select()

[void]$ws3.Activate()
[void]$ws3.Range("I7").Select()

# Leave the workbook on the first sheet, matching tabSelected="1" there.
[void]$sheet1.Activate()

Write-Host "done"
